$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 89, shifting existing rows 89:134 down to 90:135
$ws.Rows(89).Insert()

# Populate the newly inserted row 89 with the new data record
$ws.Range("A89").Value = 4
$ws.Range("B89").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C89").Value = "Los Lagos"
$ws.Range("D89").Value = 44488
$ws.Range("E89").Value = 10
$ws.Range("F89").Value = 100112032
$ws.Range("G89").Value = "Zapallo italiano"
$ws.Range("H89").Value = "Sin especificar"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 240
$ws.Range("K89").Value = 14000
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = 14500
$ws.Range("N89").Value = "`$/caja 50 unidades"
$ws.Range("O89").Value = "Región de Arica y Parinacota"
$ws.Range("P89").Value = 290
$ws.Range("Q89").Value = 50
$ws.Range("R89").Value = "Hortaliza"
